$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "6.08") are preserved exactly as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.835.10"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3
$ws.Range("D3").Value = "3.397.73"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "571.43"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("D6").Value = "162.10"
$ws.Range("E6").Value = "  +2.40%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "3.397.93"
$ws.Range("E8").Value = "  -1.22%  "

# Row 9
$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  -3.48%  "

# Row 10
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  +1.34%  "

# Row 11
$ws.Range("E11").Value = "  -1.17%  "

# Row 12
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -3.83%  "

# Row 13
$ws.Range("D13").Value = "3.977.48"
$ws.Range("E13").Value = "  -1.32%  "

# Row 14
$ws.Range("E14").Value = "  +0.91%  "

# Row 15
$ws.Range("D15").Value = "27.00"
$ws.Range("E15").Value = "  -1.85%  "

# Row 16
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -0.58%  "

# Row 17
$ws.Range("D17").Value = "63.847.17"
$ws.Range("E17").Value = "  -1.41%  "

# Row 18
$ws.Range("D18").Value = "3.413.09"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  -1.33%  "

# Row 20
$ws.Range("D20").Value = "13.58"
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("D21").Value = "376.18"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
$ws.Range("D22").Value = "7.77"
$ws.Range("E22").Value = "  -2.10%  "

# Row 23
$ws.Range("E23").Value = "  -0.18%  "

# Row 24
$ws.Range("D24").Value = "70.17"
$ws.Range("E24").Value = "  -2.29%  "

# Row 25
$ws.Range("D25").Value = "0.514"
$ws.Range("E25").Value = "  -4.88%  "

# Row 26
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  -2.92%  "

# Row 27
$ws.Range("E27").Value = "  -3.63%  "

# Row 28
$ws.Range("E28").Value = "  +0.20%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").Value = "  +0.53%  "

# Row 31
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  -4.42%  "

# Row 32
$ws.Range("D32").Value = "2.00"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33
$ws.Range("D33").Value = "22.87"
$ws.Range("E33").Value = "  -0.75%  "

# Row 34
$ws.Range("D34").Value = "7.08"
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("E35").Value = "  -4.54%  "

# Row 36
$ws.Range("D36").Value = "159.42"
$ws.Range("E36").Value = "  -0.83%  "

# Row 37
$ws.Range("D37").Value = "0.859"
$ws.Range("E37").Value = "  +10.35%  "

# Row 38
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  -2.56%  "

# Row 39
$ws.Range("D39").Value = "0.0721"
$ws.Range("E39").Value = "  -3.06%  "

# Row 40
$ws.Range("D40").Value = "42.82"
$ws.Range("E40").Value = "  -0.18%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "25.65"
$ws.Range("E41").Value = "  -1.60%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "6.46"
$ws.Range("E42").Value = "  -3.02%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.738.67"
$ws.Range("E43").Value = "  -5.47%  "

# Row 44
$ws.Range("D44").Value = "26.00"
$ws.Range("E44").Value = "  +0.61%  "

# Row 45
$ws.Range("D45").Value = "4.38"
$ws.Range("E45").Value = "  -2.70%  "

# Row 46
$ws.Range("D46").Value = "0.0306"
$ws.Range("E46").Value = "  -0.93%  "

# Row 47
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +6.88%  "

# Row 48
$ws.Range("D48").Value = "326.70"
$ws.Range("E48").Value = "  +2.92%  "

# Row 49
$ws.Range("D49").Value = "1.04"
$ws.Range("E49").Value = "  -3.83%  "

# Row 50
$ws.Range("D50").Value = "6.28"
$ws.Range("E50").Value = "  -2.65%  "

# Row 51
$ws.Range("E51").Value = "  -1.76%  "
